$wb = $excel.ActiveWorkbook

# --- Insert a new "TRANSFER" worksheet between ACCOUNT and REGULAR_TRANSFER ---
$wsAccount = $wb.Worksheets.Item("ACCOUNT")
$wsTransfer = $wb.Worksheets.Add($null, $wsAccount)
$wsTransfer.Name = "TRANSFER"

$wsTransfer.Range("A1").Value = "ID"
$wsTransfer.Range("B1").Value = "TITLE"
$wsTransfer.Range("C1").Value = "USER_ID"

$wsTransfer.Range("A2").Value = 1
$wsTransfer.Range("B2").Value = "2023/06"
$wsTransfer.Range("C2").Value = 1

$wsTransfer.Range("A3").Value = 2
$wsTransfer.Range("B3").Value = "2023/07"
$wsTransfer.Range("C3").Value = 1

$wsTransfer.Range("A4").Value = 3
$wsTransfer.Range("B4").Value = "2022/09"
$wsTransfer.Range("C4").Value = 2

$wsTransfer.Range("A1:C4").Select()

# --- Add TRANSFER_ID column (G) to TEMPORARY_TRANSFER sheet ---
$wsTemp = $wb.Worksheets.Item("TEMPORARY_TRANSFER")
$wsTemp.Range("G1").Value = "TRANSFER_ID"
$wsTemp.Range("G2").Value = 1
$wsTemp.Range("G3").Value = 2
$wsTemp.Range("G4").Value = 3

$wsTemp.Range("I14").Select()
